$wb = $excel.ActiveWorkbook

# Rename the "SCart*" sheets to their new "Cart*" names
$wb.Worksheets.Item("SCart").Name    = "Cart"
$wb.Worksheets.Item("SCartAdd").Name = "CartAdd"
$wb.Worksheets.Item("SCartUpd").Name = "CartUpdate"
$wb.Worksheets.Item("SCartDel").Name = "CartDelete"

# Fill in a previously-blank test-data cell on the CartAdd sheet with a
# whitespace-only value (3 spaces)
$ws = $wb.Worksheets.Item("CartAdd")
$ws.Range("A3").Value = "   "

# Move the active/selected tab from ChangePass to the renamed CartDelete sheet
$wb.Worksheets.Item("CartDelete").Activate()
